$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 values to the new "custom accuracy" (2 decimal place) readings.
$ws.Range("B5").Value = 20.18
$ws.Range("C5").Value = 14.82
$ws.Range("D5").Value = 1.22
$ws.Range("E5").Value = 43.86
$ws.Range("F5").Value = 35.71
$ws.Range("G5").Value = 15.85
$ws.Range("H5").Value = 56.86
$ws.Range("I5").Value = 24.43
$ws.Range("J5").Value = 10.8
$ws.Range("K5").Value = 15.98
$ws.Range("L5").Value = 17.6
$ws.Range("M5").Value = 18.54
$ws.Range("N5").Value = 5.07
$ws.Range("O5").Value = 15.79
$ws.Range("P5").Value = 22.41
$ws.Range("Q5").Value = 13.36
$ws.Range("R5").Value = 0.82
$ws.Range("S5").Value = 0.83
$ws.Range("T5").Value = 233.01
$ws.Range("U5").Value = 44.01
$ws.Range("V5").Value = 14.58
$ws.Range("W5").Value = 29.52
$ws.Range("X5").Value = 15.5
$ws.Range("Y5").Value = 2.41
$ws.Range("Z5").Value = 28.19
$ws.Range("AA5").Value = 12.87
$ws.Range("AB5").Value = 11.43
$ws.Range("AC5").Value = 13.45
$ws.Range("AD5").Value = 18.48
$ws.Range("AE5").Value = 0.55
$ws.Range("AF5").Value = 51.41
$ws.Range("AG5").Value = 8.17
$ws.Range("AH5").Value = 18.22

# Remove the last data row (row 6) entirely — the dataset now only has rows 1-5.
$ws.Rows(6).Delete()
